# Automatische test-sync: 2025-06-19 21:17:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 8 to the "Logs" sheet with the new incoming mail entry.
$logs.Range("A8").Value = "Vragen over samenwerking"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D8").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F8").Value = "2025-06-19 21:17:14"
$logs.Range("G8").Value = "Nee"

# Update the aggregated count on the "Dashboard" sheet for this category.
$dashboard.Range("B2").Value = 2

# Extend the conditional formatting ranges to cover the newly added row.
$catFcs = $logs.Range("D2:D7").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D8"))
}

$answeredFcs = $logs.Range("G2:G7").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G8"))
}
